$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "2025-10-22 03:57:32"
$ws.Range("B3").Value = "'2025-10-20"
$ws.Range("C3").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice20102025.pdf"
$ws.Range("D3").Value = "/home/runner/work/rashtriyametal_downloader/rashtriyametal_downloader/data/RashtriyaMetal/PDFs/ListPrice20102025.pdf"
